# ---------------------------------------------------------------------------
# Applies the OOXML diff to before.docx:
#   1. FirstParagraph after "Introdução" gets a new leading sentence before
#      the existing "Texto de exemplo".
#   2. First "Corpodetexto" paragraph ("Texto de exemplo") becomes a
#      citation-filled sentence.
#   3. A new "Referências" heading (with bookmark) and two ABNT-style
#      "Bibliografia" paragraphs (journal name in bold) are appended at the
#      end of the document body.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. FirstParagraph right after the "Introdução" heading -----------------
# That is the paragraph whose whole text is exactly "Texto de exemplo" and
# whose style is "FirstParagraph" (paragraph index 6).
$pIntro = $d.Paragraphs.Item(6)
$pIntro.Range.Find.Execute("Texto de exemplo", $true, $false, $false, $false, `
    $false, $true, 1, $false, `
    "Aqui será iniciada a introdução Texto de exemplo", 2) | Out-Null

# --- 2. First "Corpodetexto" paragraph --------------------------------------
# That is the first paragraph using the "Corpodetexto" style, also holding
# exactly "Texto de exemplo" (paragraph index 7).
$pCorpo = $d.Paragraphs.Item(7)
$novoTexto = "Texto de exemplo conforme (ANGELINI; HEUVELINK; KEMPEN, 2017)" + `
    " ou ANGELINI; HEUVELINK; KEMPEN (2017) e KHOSHGOFTARMANESH et al. (2018)"
$pCorpo.Range.Find.Execute("Texto de exemplo", $true, $false, $false, $false, `
    $false, $true, 1, $false, $novoTexto, 2) | Out-Null

# --- 3. Append "Referências" heading + two bibliography entries ------------

# Heading paragraph (style "Ttulo1") with a bookmark named "referencias"
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$headPara = $d.Paragraphs.Last
$headPara.Style = "Ttulo1"
$ins = $d.Range($headPara.Range.Start, $headPara.Range.Start)
$headStart = $ins.Start
$ins.InsertAfter("Referências")
$d.Bookmarks.Add("referencias", $d.Range($headStart, $headStart)) | Out-Null

# First reference paragraph (style "Bibliografia")
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Style = "Bibliografia"
$ins = $d.Range($p1.Range.Start, $p1.Range.Start)
$ins.InsertAfter("ANGELINI, M. E.; HEUVELINK, G. B. M.; KEMPEN, B. Multivariate mapping of soil with structural equation modelling.")
$ins.Collapse(0)
$ins.InsertAfter(" ")
$ins.Collapse(0)
$boldStart = $ins.Start
$ins.InsertAfter("European Journal of Soil Science")
$boldEnd = $ins.End
$d.Range($boldStart, $boldEnd).Bold = 1
$ins.Collapse(0)
$ins.InsertAfter(", v. 68, n. 5, p. 575–591, 2017.")
$ins.Collapse(0)
$ins.InsertAfter(" ")

# Second reference paragraph (style "Bibliografia")
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Style = "Bibliografia"
$ins2 = $d.Range($p2.Range.Start, $p2.Range.Start)
$ins2.InsertAfter("KHOSHGOFTARMANESH, A. H. et al. Fractionation and bioavailability of zinc (Zn) in the rhizosphere of two wheat cultivars with different Zn deficiency tolerance.")
$ins2.Collapse(0)
$ins2.InsertAfter(" ")
$ins2.Collapse(0)
$boldStart2 = $ins2.Start
$ins2.InsertAfter("Geoderma")
$boldEnd2 = $ins2.End
$d.Range($boldStart2, $boldEnd2).Bold = 1
$ins2.Collapse(0)
$ins2.InsertAfter(", v. 309, n. Supplement C, p. 1–6, 2018.")
$ins2.Collapse(0)
$ins2.InsertAfter(" ")

Write-Host "edit.ps1 applied"
